# Actualización automática 2025-09-30 15:30:09
# Applies the diff's cell value updates across the three worksheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("M10").Value = 2272.64
$ws1.Range("D11").Value = 190.08
$ws1.Range("H11").Value = 1698.3
$ws1.Range("G13").Value = 166.43
$ws1.Range("M13").Value = 13047.51
$ws1.Range("L17").Value = 525.2
$ws1.Range("D22").Value = 1419.8
$ws1.Range("H22").Value = 2293.19
$ws1.Range("I22").Value = 255.6

$ws1.Range("D23").Value = "3 de 21"
$ws1.Range("G23").Value = "1 de 21"
$ws1.Range("H23").Value = "4 de 21"
$ws1.Range("I23").Value = "2 de 21"
$ws1.Range("L23").Value = "2 de 21"
$ws1.Range("M23").Value = "9 de 21"

# ---------------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F10").Value = 2272.64
$ws2.Range("F11").Value = 6897.82
$ws2.Range("F13").Value = 14573.17
$ws2.Range("F17").Value = 6213.78
$ws2.Range("F22").Value = 10848.32
$ws2.Range("F23").Value = 61624.43

# ---------------------------------------------------------------------------
# Sheet 3: "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column F width goes from 25 -> 24 raw OOXML units. The COM layer's
# ColumnWidth property is offset from the raw <col width="..."> value by a
# constant ~5/6 (0.8333...), so compensate to land exactly on 24.
$ws3.Columns.Item(6).ColumnWidth = 23.16666666666667

$ws3.Range("D3").Value = 2983.64
$ws3.Range("E3").Value = 2520.97890386263
$ws3.Range("F3").Value = 0.5420248071862629

$ws3.Range("D5").Value = 166.43
$ws3.Range("E5").Value = -16.43000000000001
$ws3.Range("F5").Value = 1.109533333333333

$ws3.Range("D6").Value = 5690.69
$ws3.Range("E6").Value = -2783.10631853974
$ws3.Range("F6").Value = 1.957188725568165

$ws3.Range("D7").Value = 383.4
$ws3.Range("E7").Value = 503.311016287574
$ws3.Range("F7").Value = 0.4323843878755392

$ws3.Range("D11").Value = 1917.23
$ws3.Range("E11").Value = 3927.21916370549
$ws3.Range("F11").Value = 0.3280428910060774

$ws3.Range("D12").Value = 48945.04
$ws3.Range("E12").Value = -12121.3969078829
$ws3.Range("F12").Value = 1.329174299174048

$ws3.Range("D15").Value = 61181.29
$ws3.Range("E15").Value = -5756.546833866223
$ws3.Range("F15").Value = 1.103862399806007
